# Update '想去人数' (attendance interest count) figures in the F column
# across the 展览 / 演出 / 全部类型 sheets, per upstream data refresh.
$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5959
$ws.Range("F5").Value = 5959
$ws.Range("F7").Value = 3009
$ws.Range("F8").Value = 1287
$ws.Range("F12").Value = 26
$ws.Range("F13").Value = 308
$ws.Range("F14").Value = 4423
$ws.Range("F15").Value = 4423
$ws.Range("F16").Value = 101
$ws.Range("F17").Value = 92
$ws.Range("F18").Value = 128
$ws.Range("F20").Value = 195
$ws.Range("F22").Value = 6838
$ws.Range("F23").Value = 6838
$ws.Range("F24").Value = 236
$ws.Range("F26").Value = 469
$ws.Range("F27").Value = 1267
$ws.Range("F28").Value = 6265
$ws.Range("F29").Value = 1644
$ws.Range("F31").Value = 1966
$ws.Range("F32").Value = 6022
$ws.Range("F33").Value = 117
$ws.Range("F36").Value = 86
$ws.Range("F37").Value = 432
$ws.Range("F38").Value = 5928
$ws.Range("F40").Value = 190
$ws.Range("F43").Value = 8
$ws.Range("F49").Value = 349
$ws.Range("F50").Value = 2078
$ws.Range("F51").Value = 17
$ws.Range("F52").Value = 1030

# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 200
$ws.Range("F10").Value = 5

# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 5959
$ws.Range("F5").Value = 5959
$ws.Range("F7").Value = 3009
$ws.Range("F8").Value = 1287
$ws.Range("F12").Value = 200
$ws.Range("F13").Value = 308
$ws.Range("F14").Value = 4423
$ws.Range("F15").Value = 4423
$ws.Range("F16").Value = 101
$ws.Range("F17").Value = 92
$ws.Range("F18").Value = 128
$ws.Range("F20").Value = 195
$ws.Range("F22").Value = 6838
$ws.Range("F23").Value = 6838
$ws.Range("F24").Value = 236
$ws.Range("F26").Value = 469
$ws.Range("F27").Value = 1267
$ws.Range("F29").Value = 6265
$ws.Range("F30").Value = 1644
$ws.Range("F33").Value = 1966
$ws.Range("F34").Value = 6022
$ws.Range("F35").Value = 117
$ws.Range("F36").Value = 5
$ws.Range("F38").Value = 86
$ws.Range("F39").Value = 432
$ws.Range("F40").Value = 5928
$ws.Range("F42").Value = 190
$ws.Range("F45").Value = 8
$ws.Range("F50").Value = 349
$ws.Range("F51").Value = 17

